# Append a new data row (row 78) to each of the 4 worksheets, following the
# existing pattern of rows in the log/database sheets.

$wb = $excel.ActiveWorkbook

$dateValue = 45864.43591435185
$dateFormat = "YYYY-MM-DD HH:MM:SS"

$rowsData = @(
    @{
        Sheet = 1
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x3C"
        E = "0x14"
        F = 380
        G = "7.598631275147109e+23" -as [double]
        H = 316
        I = 14
    },
    @{
        Sheet = 2
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x40"
        E = "0xe"
        F = 380
        G = "5.68432987514711e+23" -as [double]
        H = 320
        I = 14
    },
    @{
        Sheet = 3
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x78"
        E = "0x7"
        F = 130
        G = "5.68631262647114e+23" -as [double]
        H = 120
        I = 7
    },
    @{
        Sheet = 4
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x78"
        E = "0x3"
        F = 130
        G = "9.85046333984776e+23" -as [double]
        H = 120
        I = 3
    }
)

foreach ($rowData in $rowsData) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)

    $ws.Range("A78").Value = $dateValue
    $ws.Range("A78").NumberFormat = $dateFormat

    $ws.Range("B78").Value = $rowData.B
    $ws.Range("C78").Value = $rowData.C
    $ws.Range("D78").Value = $rowData.D
    $ws.Range("E78").Value = $rowData.E
    $ws.Range("F78").Value = $rowData.F
    $ws.Range("G78").Value = $rowData.G
    $ws.Range("H78").Value = $rowData.H
    $ws.Range("I78").Value = $rowData.I
}
